$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused column I (was all blank/style-only cells)
$ws.Columns.Item(9).Delete()

# Update recalculated values in the existing "CRP_nr" row (row 24) and
# the "pasture_nr" row (row 25)
$ws.Range("H24").Value = 56.813400268554688
$ws.Range("H25").Value = 16.878217697143555

# Add the new "range_nr" row (row 26), copying formatting from row 25
$ws.Range("A25:H25").Copy()
$ws.Range("A26").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A26").Value = "range_nr"
$ws.Range("E26").Value = 18.795345306396484
$ws.Range("F26").Value = 17.547082901000977
$ws.Range("G26").Value = 18.402990341186523
$ws.Range("H26").Value = 16.878217697143555
